$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.869.62"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.887.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.62%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7578"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.87%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3121"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.40"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.98%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07147"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.48%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08533"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.94%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7609"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.924.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.358"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.95%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.89%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.133"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.77%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.900.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007809"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.147.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.84%  "

$ws.Range("E24").Value = "  +0.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1609"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.374"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.29%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.43%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.025"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.90%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.520"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.532"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.98%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.479"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.51%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.106"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05418"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.71%  "

$ws.Range("E35").Value = "  -1.26%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7437"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9996"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.709"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01940"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.776"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4455"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.103.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.17%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.075"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.59%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8598"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.867"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.09%  "

$ws.Range("E49").Value = "  +1.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.052"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.043.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.35%  "
